$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (F12.bet) - replace N/A placeholders with actual values
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "98.9%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "98.4%"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "99.8%"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "9.5"

# Row 4 (Luva bet) - update Reclamacoes Respondidas
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "98.1%"
